$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.139.70'
$ws.Range('E2').Value = '  -0.18%  '
$ws.Range('D3').Value = '3.553.39'
$ws.Range('E3').Value = '  +1.65%  '
$ws.Range('D5').Value = '''605.24'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.05%  '
$ws.Range('D6').Value = '''144.06'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -0.28%  '
$ws.Range('D7').Value = '3.552.16'
$ws.Range('E7').Value = '  +1.66%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('E9').Value = '  +3.01%  '
$ws.Range('E10').Value = '  -0.34%  '
$ws.Range('D11').Value = '''7.80'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -3.42%  '
$ws.Range('E12').Value = '  -0.33%  '
$ws.Range('E13').Value = '  +1.66%  '
$ws.Range('E14').Value = '  +0.66%  '
$ws.Range('E15').Value = '  -1.22%  '
$ws.Range('D16').Value = '3.565.15'
$ws.Range('E16').Value = '  +2.19%  '
$ws.Range('D17').Value = '66.190.28'
$ws.Range('E17').Value = '  -0.15%  '
$ws.Range('E18').Value = '  -1.08%  '
$ws.Range('D19').Value = '''11.34'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +5.51%  '
$ws.Range('E20').Value = '  +0.36%  '
$ws.Range('D21').Value = '''14.66'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -1.36%  '
$ws.Range('E22').Value = '  +0.81%  '
$ws.Range('E23').Value = '  +2.28%  '
$ws.Range('D24').Value = '''79.93'
$ws.Range('D24').Style = "Normal"
$ws.Range('D25').Value = '3.695.43'
$ws.Range('E25').Value = '  +1.78%  '
$ws.Range('E26').Value = '  -0.05%  '
$ws.Range('E27').Value = '  -0.62%  '
$ws.Range('E28').Value = '  +0.54%  '
$ws.Range('D29').Value = '''9.07'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -2.39%  '
$ws.Range('D30').Value = '''7.83'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -1.25%  '
$ws.Range('D31').Value = '''1.00'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +0.20%  '
$ws.Range('D32').Value = '3.550.00'
$ws.Range('E32').Value = '  +1.96%  '
$ws.Range('D33').Value = '''25.40'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +1.28%  '
$ws.Range('E34').Value = '  -1.91%  '
$ws.Range('E35').Value = '  -9.13%  '
$ws.Range('D37').Value = '''7.79'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +0.65%  '
$ws.Range('E38').Value = '  -1.36%  '
$ws.Range('E39').Value = '  -1.65%  '
$ws.Range('D40').Value = '''173.72'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +2.12%  '
$ws.Range('D41').Value = '''0.0843'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -1.96%  '
$ws.Range('D42').Value = '''5.17'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +0.01%  '
$ws.Range('D43').Value = '''0.889'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +0.84%  '
$ws.Range('D44').Value = '''1.91'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +0.58%  '
$ws.Range('D45').Value = '''46.05'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +1.38%  '
$ws.Range('E46').Value = '  +0.03%  '
$ws.Range('E47').Value = '  -1.55%  '
$ws.Range('D48').Value = '''24.82'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -4.42%  '
$ws.Range('D49').Value = '''2.40'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -1.22%  '
$ws.Range('D50').Value = '''7.10'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -0.65%  '
$ws.Range('D51').Value = '''22.85'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +2.75%  '
